$wb = $excel.ActiveWorkbook

$compData = $wb.Worksheets.Item("Component Data")
$compData.Range("F3").Value = 700
$compData.Range("I4").Value = 1000
$compData.Range("F8").Value = 0.5
$compData.Range("G8").Value = 24
$compData.Range("H8").Value = 1
$compData.Range("F10").Value = 7.5
$compData.Range("H10").Value = ""
